$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Arduino" characteristic to F9 (same row as the
# "0.513L - Laboratorio MAC" lab entry).
$ws.Range("F9").Value = "Arduino"

# Move the active selection to F9 (matches the edit's new cursor position).
$ws.Range("F9").Select()
